$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 3352.3333
$ws.Range("J18").Value = 5623
$ws.Range("L18").Value = 5623
$ws.Range("N18").Value = -6191
# Row 33
$ws.Range("H33").Value = 263.04
$ws.Range("I33").Value = 243.05263
$ws.Range("K33").Value = 243.05263
$ws.Range("M33").Value = -14.05262999999999
# Row 34
$ws.Range("H34").Value = 1863.5555
$ws.Range("I34").Value = 1863.5555
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1863.5555
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1660.5555
$ws.Range("N34").ClearContents()
# Row 36
$ws.Range("H36").Value = 1863.5555
$ws.Range("I36").Value = 1863.5555
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1863.5555
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1148.5555
$ws.Range("N36").ClearContents()
# Row 38
$ws.Range("H38").Value = 10808.25
$ws.Range("I38").Value = 13411.333
$ws.Range("J38").Value = 2999
$ws.Range("K38").Value = 40233.999
$ws.Range("L38").Value = 8997
$ws.Range("M38").Value = -39861.999
$ws.Range("N38").Value = -9741
# Row 52
$ws.Range("H52").Value = 720
$ws.Range("I52").Value = 300
$ws.Range("J52").Value = 750
$ws.Range("K52").Value = 900
$ws.Range("L52").Value = 2250
$ws.Range("M52").Value = -740
$ws.Range("N52").Value = -2570
# Row 61
$ws.Range("H61").Value = 1034.4
$ws.Range("I61").Value = 1034.4
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3103.2
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2931.2
$ws.Range("N61").ClearContents()
# Row 86
$ws.Range("H86").Value = 5648.9375
$ws.Range("I86").Value = 4413.636
$ws.Range("K86").Value = 4413.636
$ws.Range("M86").Value = -3290.636
# Row 89
$ws.Range("H89").Value = 5648.9375
$ws.Range("I89").Value = 4413.636
$ws.Range("K89").Value = 22068.18
$ws.Range("M89").Value = -16452.18
# Row 92
$ws.Range("H92").Value = 1233.5358
$ws.Range("J92").Value = 1999.75
$ws.Range("L92").Value = 1999.75
$ws.Range("N92").Value = -4495.75
# Row 106
$ws.Range("H106").Value = 4030.182
$ws.Range("I106").Value = 4338.857
$ws.Range("K106").Value = 4338.857
$ws.Range("M106").Value = -3707.857
# Row 116
$ws.Range("H116").Value = 7999.1113
$ws.Range("J116").Value = 7537.769
$ws.Range("L116").Value = 7537.769
$ws.Range("N116").Value = -14421.769
# Row 138
$ws.Range("H138").Value = 1920136.1
$ws.Range("I138").Value = 5690.28
$ws.Range("J138").Value = 3370473.8
$ws.Range("K138").Value = 17070.84
$ws.Range("L138").Value = 10111421.4
$ws.Range("M138").Value = -11930.84
$ws.Range("N138").Value = -10121701.4
# Row 141
$ws.Range("H141").Value = 5029.4
$ws.Range("J141").Value = 14490.833
$ws.Range("L141").Value = 43472.499
$ws.Range("N141").Value = -53832.499

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3469.3572
$ws.Range("I2").Value = 3271
$ws.Range("J2").Value = 3888.111
$ws.Range("K2").Value = 3271
$ws.Range("L2").Value = 3888.111
$ws.Range("M2").Value = -3158
$ws.Range("N2").Value = -4114.111
# Row 32
$ws.Range("H32").Value = 7877.9185
$ws.Range("I32").Value = 6308.3335
$ws.Range("J32").Value = 13999.3
$ws.Range("K32").Value = 6308.3335
$ws.Range("L32").Value = 13999.3
$ws.Range("M32").Value = -6021.3335
$ws.Range("N32").Value = -14573.3
# Row 45
$ws.Range("H45").Value = 6563.857
$ws.Range("I45").Value = 9624.75
$ws.Range("J45").Value = 2482.6667
$ws.Range("K45").Value = 9624.75
$ws.Range("L45").Value = 2482.6667
$ws.Range("M45").Value = -9247.75
$ws.Range("N45").Value = -3236.6667
# Row 92
$ws.Range("H92").Value = 60001
$ws.Range("J92").Value = 60001
$ws.Range("L92").Value = 60001
$ws.Range("N92").Value = -64993
# Row 95
$ws.Range("H95").Value = 71736
$ws.Range("J95").Value = 71736
$ws.Range("L95").Value = 71736
$ws.Range("N95").Value = -77228
# Row 110
$ws.Range("H110").Value = 942.6923
$ws.Range("I110").Value = 880.4
$ws.Range("K110").Value = 880.4
$ws.Range("M110").Value = 1164.6
# Row 116
$ws.Range("H116").Value = 3469.3572
$ws.Range("I116").Value = 3271
$ws.Range("J116").Value = 3888.111
$ws.Range("K116").Value = 3271
$ws.Range("L116").Value = 3888.111
$ws.Range("M116").Value = -977
$ws.Range("N116").Value = -8476.111000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3469.3572
$ws.Range("I3").Value = 3271
$ws.Range("J3").Value = 3888.111
$ws.Range("K3").Value = 3271
$ws.Range("L3").Value = 3888.111
$ws.Range("M3").Value = -3157
$ws.Range("N3").Value = -4116.111
# Row 86
$ws.Range("H86").Value = 617308.7
$ws.Range("J86").Value = 2714297.5
$ws.Range("L86").Value = 2714297.5
$ws.Range("N86").Value = -2716543.5
# Row 89
$ws.Range("H89").Value = 617308.7
$ws.Range("J89").Value = 2714297.5
$ws.Range("L89").Value = 13571487.5
$ws.Range("N89").Value = -13582719.5
# Row 99
$ws.Range("H99").Value = 3494
$ws.Range("I99").Value = 2812
$ws.Range("K99").Value = 2812
$ws.Range("M99").Value = -1314
# Row 105
$ws.Range("H105").Value = 6899.567
$ws.Range("I105").Value = 8657.474
$ws.Range("J105").Value = 3863.182
$ws.Range("K105").Value = 8657.474
$ws.Range("L105").Value = 3863.182
$ws.Range("M105").Value = -6910.474
$ws.Range("N105").Value = -7357.182

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 1239.8948
$ws.Range("I4").Value = 1336.8125
$ws.Range("J4").Value = 723
$ws.Range("K4").Value = 1336.8125
$ws.Range("L4").Value = 723
$ws.Range("M4").Value = -1224.8125
$ws.Range("N4").Value = -947
# Row 16
$ws.Range("H16").Value = 1748.875
$ws.Range("I16").Value = 873.25
$ws.Range("J16").Value = 2624.5
$ws.Range("K16").Value = 873.25
$ws.Range("L16").Value = 2624.5
$ws.Range("M16").Value = -586.25
$ws.Range("N16").Value = -3198.5
# Row 69
$ws.Range("H69").Value = 5890.5713
$ws.Range("I69").Value = 5890.5713
$ws.Range("K69").Value = 5890.5713
$ws.Range("M69").Value = -5141.5713
# Row 72
$ws.Range("H72").Value = 5890.5713
$ws.Range("I72").Value = 5890.5713
$ws.Range("K72").Value = 17671.7139
$ws.Range("M72").Value = -13927.7139
# Row 107
$ws.Range("H107").Value = 826.45
$ws.Range("I107").Value = 889.0909
$ws.Range("J107").Value = 749.8889
$ws.Range("K107").Value = 889.0909
$ws.Range("L107").Value = 749.8889
$ws.Range("M107").Value = 1030.9091
$ws.Range("N107").Value = -4589.8889
# Row 113
$ws.Range("H113").Value = 1748.875
$ws.Range("I113").Value = 873.25
$ws.Range("J113").Value = 2624.5
$ws.Range("K113").Value = 873.25
$ws.Range("L113").Value = 2624.5
$ws.Range("M113").Value = 1296.75
$ws.Range("N113").Value = -6964.5
# Row 134
$ws.Range("H134").Value = 14536.6045
$ws.Range("I134").Value = 14994.718
$ws.Range("K134").Value = 44984.154
$ws.Range("M134").Value = -42449.154

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 92
$ws.Range("H92").Value = 129.85715
$ws.Range("I92").Value = 126.666664
$ws.Range("K92").Value = 379.999992
$ws.Range("M92").Value = 868.000008
# Row 109
$ws.Range("H109").Value = 5261
$ws.Range("I109").Value = 616.2222
$ws.Range("K109").Value = 1848.6666
$ws.Range("M109").Value = -808.6666
# Row 122
$ws.Range("H122").Value = 1460.5
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1460.5
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 13144.5
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -18044.5
# Row 129
$ws.Range("H129").Value = 7036.875
$ws.Range("I129").Value = 2765
$ws.Range("J129").Value = 8460.833000000001
$ws.Range("K129").Value = 8295
$ws.Range("L129").Value = 25382.499
$ws.Range("M129").Value = -3295
$ws.Range("N129").Value = -35382.499
# Row 132
$ws.Range("H132").Value = 2195.7646
$ws.Range("I132").Value = 2376.5
$ws.Range("J132").Value = 2171.6667
$ws.Range("K132").Value = 21388.5
$ws.Range("L132").Value = 19545.0003
$ws.Range("M132").Value = -18858.5
$ws.Range("N132").Value = -24605.0003

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2192.0454
$ws.Range("I132").Value = 1840.125
$ws.Range("J132").Value = 3130.5
$ws.Range("K132").Value = 5520.375
$ws.Range("L132").Value = 9391.5
$ws.Range("M132").Value = -2990.375
$ws.Range("N132").Value = -14451.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 159919.73
$ws.Range("I2").Value = 164705.19
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 164705.19
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -164593.19
$ws.Range("N2").Value = -2224
# Row 7
$ws.Range("H7").Value = 6774.077
$ws.Range("I7").Value = 6616.375
$ws.Range("J7").Value = 7026.4
$ws.Range("K7").Value = 6616.375
$ws.Range("L7").Value = 7026.4
$ws.Range("M7").Value = -6504.375
$ws.Range("N7").Value = -7250.4
# Row 126
$ws.Range("H126").Value = 6774.077
$ws.Range("I126").Value = 6616.375
$ws.Range("J126").Value = 7026.4
$ws.Range("K126").Value = 19849.125
$ws.Range("L126").Value = 21079.2
$ws.Range("M126").Value = -17379.125
$ws.Range("N126").Value = -26019.2
# Row 132
$ws.Range("H132").Value = 3467.72
$ws.Range("I132").Value = 2749.6875
$ws.Range("K132").Value = 8249.0625
$ws.Range("M132").Value = -5719.0625

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 6983334
$ws.Range("I2").Value = 10350001
$ws.Range("J2").Value = 250000
$ws.Range("K2").Value = 10350001
$ws.Range("L2").Value = 250000
$ws.Range("M2").Value = -10349889
$ws.Range("N2").Value = -250224
# Row 107
$ws.Range("H107").Value = 373.63333
$ws.Range("I107").Value = 350.13635
$ws.Range("J107").Value = 438.25
$ws.Range("K107").Value = 1050.40905
$ws.Range("L107").Value = 1314.75
$ws.Range("M107").Value = 869.59095
$ws.Range("N107").Value = -5154.75
# Row 132
$ws.Range("H132").Value = 4637.864
$ws.Range("I132").Value = 3440.6
$ws.Range("K132").Value = 10321.8
$ws.Range("M132").Value = -7791.799999999999
# Row 136
$ws.Range("H136").Value = 4661.8667
$ws.Range("I136").Value = 997
$ws.Range("K136").Value = 2991
$ws.Range("M136").Value = -441

